# Updated cryptos list values (refreshed data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.049.98"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.398.55"
$ws.Range("E3").Value = "  -4.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "476.83"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("D5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.55"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("D6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.499"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D8").NumberFormat = "General"

$ws.Range("D9").Value = "2.407.78"
$ws.Range("E9").Value = "  -4.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.325"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D12").NumberFormat = "General"

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").Value = "2.830.94"
$ws.Range("E14").Value = "  -3.77%  "

$ws.Range("D15").Value = "56.366.60"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.39"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D17").NumberFormat = "General"

$ws.Range("D18").Value = "2.407.40"
$ws.Range("E18").Value = "  -4.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.46"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.81"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.78"
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("D21").NumberFormat = "General"

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.03"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D24").NumberFormat = "General"

$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.396"
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("D26").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "General"

$ws.Range("D28").Value = "2.534.91"
$ws.Range("E28").Value = "  -2.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D29").NumberFormat = "General"

$ws.Range("D30").Value = "0.0₃0774"
$ws.Range("E30").Value = "  -1.65%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.77"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "General"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.86"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "General"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D34").NumberFormat = "General"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.02"
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("D35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D36").NumberFormat = "General"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.61"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "General"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.70"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.35"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0546"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("D43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.593"
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("D44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("D45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.23"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "258.18"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.62"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0222"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "General"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.18"
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("D50").NumberFormat = "General"

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.800.81"
$ws.Range("E51").Value = "  -8.01%  "
